# Updated symbol list on Tue Jan 24 07:51:07 UTC 2023 with GitHub Actions
# Refresh cryptocurrency price / 1h volume-change figures in columns D (Price)
# and E (Volume(1h)) on Sheet1. All values are stored as literal text
# (matching the workbook's existing inlineStr cells), so each value is
# written with a leading apostrophe to stop Excel from reinterpreting the
# numeric-looking / percent-looking text as a real number, and the style is
# reset to "Normal" afterwards so no unintended text-format style gets
# attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.95%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'36.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.38%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.160"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.12%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08254"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.24%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.155"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.24%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.013"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.91%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.145"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.01%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.98%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1010"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.56%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.26%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09234"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'7.57%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03631"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.71%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09926"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001435"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.02%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005738"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.458"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.14%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.801"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'12.56%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3378"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.42%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.213"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'9.39%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.51%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2254"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.25%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04600"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.11%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001248"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.66%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004734"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.86%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-21.84%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004509"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-5.24%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02006"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'8.90%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04918"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.13%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.63%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1401"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.01%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007837"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.22%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002108"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.00%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01178"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.21%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006459"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.65%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'31.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-11.99%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.98%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.01%"
$ws.Range("E51").Style = "Normal"
